# Automatische test-sync: 2025-08-14 21:01:50
# Adds the new "Terugbetaling" log entry (row 19) to the "Logs" sheet,
# extends the conditional-formatting ranges to cover the new row, and
# bumps the "Intern verzoek / Actie voor medewerker" tally on "Dashboard".

$wb = $excel.ActiveWorkbook

# --- Logs sheet: append row 19 --------------------------------------------
$ws = $wb.Worksheets.Item("Logs")

$ws.Cells.Item(19, 1).Value = "Terugbetaling"
$ws.Cells.Item(19, 2).Value = "support@testbedrijf123.nl"
$ws.Cells.Item(19, 3).Value = "Ik heb mijn retour gestuurd maar nog geen geld terug."
$ws.Cells.Item(19, 4).Value = "Intern verzoek / Actie voor medewerker"
$ws.Cells.Item(19, 5).Value = "Bedankt, we hebben dit doorgestuurd naar retour@testbedrijf123.nl."
$ws.Cells.Item(19, 6).Value = "2025-08-14 21:01:26"
$ws.Cells.Item(19, 7).Value = "Nee"
$ws.Cells.Item(19, 8).Value = "Ja"
$ws.Cells.Item(19, 9).Value = "Nee"
$ws.Cells.Item(19, 10).Value = "Nee"

# --- Extend conditional formatting ranges from row 18 to row 19 -----------
$ws.Range("D2:D18").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D19"))
$ws.Range("G2:G18").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G19"))
$ws.Range("H2:H18").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H2:H19"))
$ws.Range("I2:I18").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("I2:I19"))
$ws.Range("J2:J18").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("J2:J19"))

# --- Dashboard sheet: bump the count for "Intern verzoek / ..." -----------
$ws2 = $wb.Worksheets.Item("Dashboard")
$ws2.Range("B2").Value = 13
